# Scheduled-runner refresh of the Sheets workbook (Ultima_Profits).
# Re-applies freshly pulled Market Board prices to the Leve profit tables:
# for every affected row we rewrite currentAveragePrice/NQ/HQ, LevePriceNQ/HQ
# and the resulting LeveProfitNQ/HQ columns (H:N). A few rows end up with no
# HQ-or-NQ profit figure at all (item not listable that way any more) -- for
# those the now-stale profit cell is cleared instead of being overwritten.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 526856.8
$ws.Range("I38").Value = 125
$ws.Range("J38").Value = 588825.25
$ws.Range("K38").Value = 375
$ws.Range("L38").Value = 1766475.75
$ws.Range("M38").Value = -3
$ws.Range("N38").Value = -1767219.75
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("H62").Value = 3400
$ws.Range("I62").Value = 5500
$ws.Range("J62").Value = 1300
$ws.Range("K62").Value = 5500
$ws.Range("L62").Value = 1300
$ws.Range("M62").Value = -4876
$ws.Range("N62").Value = -2548
$ws.Range("H65").Value = 3400
$ws.Range("I65").Value = 5500
$ws.Range("J65").Value = 1300
$ws.Range("K65").Value = 27500
$ws.Range("L65").Value = 6500
$ws.Range("M65").Value = -24380
$ws.Range("N65").Value = -12740
$ws.Range("H137").Value = 11112673
$ws.Range("I137").Value = 1313.8572
$ws.Range("J137").Value = 18183538
$ws.Range("K137").Value = 3941.5716
$ws.Range("L137").Value = 54550614
$ws.Range("M137").Value = -1391.5716
$ws.Range("N137").Value = -54555714
$ws.Range("H138").Value = 6098757
$ws.Range("I138").Value = 977.65515
$ws.Range("J138").Value = 20835056
$ws.Range("K138").Value = 2932.96545
$ws.Range("L138").Value = 62505168
$ws.Range("M138").Value = 2207.03455
$ws.Range("N138").Value = -62515448
$ws.Range("H141").Value = 2821.5
$ws.Range("I141").Value = 2066.4285
$ws.Range("J141").Value = 4583.3335
$ws.Range("K141").Value = 6199.2855
$ws.Range("L141").Value = 13750.0005
$ws.Range("M141").Value = -1019.2855
$ws.Range("N141").Value = -24110.0005
$ws.Range("M47").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 183
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 199.5
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 199.5
$ws.Range("M4").Value = -34
$ws.Range("N4").Value = -431.5
$ws.Range("H61").Value = 12822635
$ws.Range("I61").Value = 12822635
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 12822635
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -12822423
$ws.Range("H110").Value = 1132.5
$ws.Range("I110").Value = 1132.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1132.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 912.5
$ws.Range("H122").Value = 2094.6562
$ws.Range("I122").Value = 2144.44
$ws.Range("K122").Value = 6433.32
$ws.Range("M122").Value = -3983.32
$ws.Range("H136").Value = 12822635
$ws.Range("I136").Value = 12822635
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 38467905
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -38465355
$ws.Range("N61").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2443.125
$ws.Range("I107").Value = 2474.7144
$ws.Range("K107").Value = 2474.7144
$ws.Range("M107").Value = -554.7143999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1337.5
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1675
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1675
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -2249
$ws.Range("H19").Value = 115.833336
$ws.Range("I19").Value = 115.833336
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 115.833336
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 54.166664
$ws.Range("H24").Value = 115.833336
$ws.Range("I24").Value = 115.833336
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 115.833336
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 54.166664
$ws.Range("H58").Value = 2384.9546
$ws.Range("I58").Value = 1165.625
$ws.Range("K58").Value = 1165.625
$ws.Range("M58").Value = -962.625
$ws.Range("H113").Value = 1337.5
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1675
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1675
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -6015
$ws.Range("H132").Value = 4517.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4517.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13552.5
$ws.Range("N132").Value = -18612.5
$ws.Range("H134").Value = 612194.7
$ws.Range("I134").Value = 1446
$ws.Range("K134").Value = 4338
$ws.Range("M134").Value = -1803
$ws.Range("H136").Value = 2384.9546
$ws.Range("I136").Value = 1165.625
$ws.Range("K136").Value = 3496.875
$ws.Range("M136").Value = -946.875
$ws.Range("N19").ClearContents()
$ws.Range("N24").ClearContents()
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 613.3333
$ws.Range("I5").Value = 408.57144
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 1225.71432
$ws.Range("L5").Value = 2700
$ws.Range("M5").Value = -1113.71432
$ws.Range("N5").Value = -2924
$ws.Range("H63").Value = 6666.1665
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 7599.4
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 22798.2
$ws.Range("M63").Value = -5251
$ws.Range("N63").Value = -24296.2
$ws.Range("H66").Value = 6666.1665
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 7599.4
$ws.Range("K66").Value = 18000
$ws.Range("L66").Value = 68394.59999999999
$ws.Range("M66").Value = -14256
$ws.Range("N66").Value = -75882.59999999999
$ws.Range("H118").Value = 3253.7693
$ws.Range("J118").Value = 3580.25
$ws.Range("L118").Value = 10740.75
$ws.Range("N118").Value = -13226.75
$ws.Range("H132").Value = 841
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 841
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 7569
$ws.Range("N132").Value = -12629
$ws.Range("H135").Value = 613.3333
$ws.Range("I135").Value = 408.57144
$ws.Range("J135").Value = 900
$ws.Range("K135").Value = 3677.14296
$ws.Range("L135").Value = 8100
$ws.Range("M135").Value = -1142.14296
$ws.Range("N135").Value = -13170
$ws.Range("H138").Value = 3037.3
$ws.Range("I138").Value = 1955
$ws.Range("J138").Value = 7366.5
$ws.Range("K138").Value = 5865
$ws.Range("L138").Value = 22099.5
$ws.Range("M138").Value = -725
$ws.Range("N138").Value = -32379.5
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 18999.5
$ws.Range("J111").Value = 18999.5
$ws.Range("L111").Value = 18999.5
$ws.Range("N111").Value = -25133.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 2000
$ws.Range("J24").Value = 2000
$ws.Range("L24").Value = 2000
$ws.Range("N24").Value = -2686
$ws.Range("H40").Value = 6000.4375
$ws.Range("I40").Value = 10840
$ws.Range("K40").Value = 10840
$ws.Range("M40").Value = -10704
$ws.Range("H132").Value = 7582290
$ws.Range("I132").Value = 4611.0234
$ws.Range("J132").Value = 21749256
$ws.Range("K132").Value = 13833.0702
$ws.Range("L132").Value = 65247768
$ws.Range("M132").Value = -11303.0702
$ws.Range("N132").Value = -65252828
$ws.Range("H134").Value = 56052.723
$ws.Range("J134").Value = 56052.723
$ws.Range("L134").Value = 56052.723
$ws.Range("N134").Value = -66192.723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H100").Value = 2815
$ws.Range("I100").Value = 3830.8333
$ws.Range("J100").Value = 1596
$ws.Range("K100").Value = 7661.6666
$ws.Range("L100").Value = 3192
$ws.Range("M100").Value = -7120.6666
$ws.Range("N100").Value = -4274
$ws.Range("H132").Value = 1348.1111
$ws.Range("I132").Value = 1108.4231
$ws.Range("J132").Value = 1971.3
$ws.Range("K132").Value = 3325.2693
$ws.Range("L132").Value = 5913.9
$ws.Range("M132").Value = -795.2692999999999
$ws.Range("N132").Value = -10973.9
$ws.Range("H136").Value = 1285.3182
$ws.Range("I136").Value = 965.5714
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 2896.7142
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -346.7142000000003
$ws.Range("N136").Value = -29100
$ws.Range("N16").ClearContents()

